$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 1133.3334
$ws.Cells.Item(4, 9).Value = 940
$ws.Cells.Item(4, 11).Value = 940
$ws.Cells.Item(4, 13).Value = -826

$ws.Cells.Item(29, 8).Value = 624.75
$ws.Cells.Item(29, 9).Value = 499.66666
$ws.Cells.Item(29, 10).Value = 1000
$ws.Cells.Item(29, 11).Value = 1498.99998
$ws.Cells.Item(29, 12).Value = 3000
$ws.Cells.Item(29, 13).Value = -1217.99998
$ws.Cells.Item(29, 14).Value = -3562

$ws.Cells.Item(98, 8).Value = 35464.22
$ws.Cells.Item(98, 9).Value = 1383.75
$ws.Cells.Item(98, 10).Value = 72642.91
$ws.Cells.Item(98, 11).Value = 1383.75
$ws.Cells.Item(98, 12).Value = 72642.91
$ws.Cells.Item(98, 13).Value = 114.25
$ws.Cells.Item(98, 14).Value = -75638.91

$ws.Cells.Item(120, 8).Value = 43000
$ws.Cells.Item(120, 10).Value = 43000
$ws.Cells.Item(120, 12).Value = 43000
$ws.Cells.Item(120, 14).Value = -52676

$ws.Cells.Item(122, 8).Value = 35464.22
$ws.Cells.Item(122, 9).Value = 1383.75
$ws.Cells.Item(122, 10).Value = 72642.91
$ws.Cells.Item(122, 11).Value = 4151.25
$ws.Cells.Item(122, 12).Value = 217928.73
$ws.Cells.Item(122, 13).Value = -1701.25
$ws.Cells.Item(122, 14).Value = -222828.73

$ws.Cells.Item(137, 8).Value = 2490150.2
$ws.Cells.Item(137, 9).Value = 4816853
$ws.Cells.Item(137, 10).Value = 8333.666999999999
$ws.Cells.Item(137, 11).Value = 14450559
$ws.Cells.Item(137, 12).Value = 25001.001
$ws.Cells.Item(137, 13).Value = -14448009
$ws.Cells.Item(137, 14).Value = -30101.001

$ws.Cells.Item(138, 8).Value = 2719.6316
$ws.Cells.Item(138, 9).Value = 1787.7916
$ws.Cells.Item(138, 10).Value = 3397.3333
$ws.Cells.Item(138, 11).Value = 5363.3748
$ws.Cells.Item(138, 12).Value = 10191.9999
$ws.Cells.Item(138, 13).Value = -223.3747999999996
$ws.Cells.Item(138, 14).Value = -20471.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 20643.572
$ws.Cells.Item(32, 9).Value = 20996.322
$ws.Cells.Item(32, 10).Value = 19232.572
$ws.Cells.Item(32, 11).Value = 20996.322
$ws.Cells.Item(32, 12).Value = 19232.572
$ws.Cells.Item(32, 13).Value = -20709.322
$ws.Cells.Item(32, 14).Value = -19806.572

$ws.Cells.Item(45, 8).Value = 2070.389
$ws.Cells.Item(45, 9).Value = 1768.0714
$ws.Cells.Item(45, 10).Value = 3128.5
$ws.Cells.Item(45, 11).Value = 1768.0714
$ws.Cells.Item(45, 12).Value = 3128.5
$ws.Cells.Item(45, 13).Value = -1391.0714
$ws.Cells.Item(45, 14).Value = -3882.5

$ws.Cells.Item(61, 8).Value = 2497.6667
$ws.Cells.Item(61, 9).Value = 1675.125
$ws.Cells.Item(61, 10).Value = 2908.9375
$ws.Cells.Item(61, 11).Value = 1675.125
$ws.Cells.Item(61, 12).Value = 2908.9375
$ws.Cells.Item(61, 13).Value = -1463.125
$ws.Cells.Item(61, 14).Value = -3332.9375

$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 14).ClearContents()

$ws.Cells.Item(122, 8).Value = 1923.5294
$ws.Cells.Item(122, 9).Value = 1565.5
$ws.Cells.Item(122, 11).Value = 4696.5
$ws.Cells.Item(122, 13).Value = -2246.5

$ws.Cells.Item(136, 8).Value = 2497.6667
$ws.Cells.Item(136, 9).Value = 1675.125
$ws.Cells.Item(136, 10).Value = 2908.9375
$ws.Cells.Item(136, 11).Value = 5025.375
$ws.Cells.Item(136, 12).Value = 8726.8125
$ws.Cells.Item(136, 13).Value = -2475.375
$ws.Cells.Item(136, 14).Value = -13826.8125

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(115, 8).Value = 19698.545
$ws.Cells.Item(115, 10).Value = 19698.545
$ws.Cells.Item(115, 12).Value = 19698.545
$ws.Cells.Item(115, 14).Value = -22832.545

$ws.Cells.Item(134, 8).Value = 3997.7273
$ws.Cells.Item(134, 9).Value = 3593.0715
$ws.Cells.Item(134, 10).Value = 4705.875
$ws.Cells.Item(134, 11).Value = 10779.2145
$ws.Cells.Item(134, 12).Value = 14117.625
$ws.Cells.Item(134, 13).Value = -8244.2145
$ws.Cells.Item(134, 14).Value = -19187.625

$ws.Cells.Item(135, 8).Value = 45000
$ws.Cells.Item(135, 10).Value = 45000
$ws.Cells.Item(135, 12).Value = 45000
$ws.Cells.Item(135, 14).Value = -55140

$ws.Cells.Item(137, 8).Value = 60114.6
$ws.Cells.Item(137, 10).Value = 60114.6
$ws.Cells.Item(137, 12).Value = 60114.6
$ws.Cells.Item(137, 14).Value = -70314.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3549370.2
$ws.Cells.Item(31, 9).Value = 1615.6097
$ws.Cells.Item(31, 10).Value = 6293860
$ws.Cells.Item(31, 11).Value = 1615.6097
$ws.Cells.Item(31, 12).Value = 6293860
$ws.Cells.Item(31, 13).Value = -1320.6097
$ws.Cells.Item(31, 14).Value = -6294450

$ws.Cells.Item(34, 8).Value = 3549370.2
$ws.Cells.Item(34, 9).Value = 1615.6097
$ws.Cells.Item(34, 10).Value = 6293860
$ws.Cells.Item(34, 11).Value = 1615.6097
$ws.Cells.Item(34, 12).Value = 6293860
$ws.Cells.Item(34, 13).Value = -1413.6097
$ws.Cells.Item(34, 14).Value = -6294264

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value = 43588.57
$ws.Cells.Item(18, 9).Value = 50803.332
$ws.Cells.Item(18, 10).Value = 300
$ws.Cells.Item(18, 11).Value = 152409.996
$ws.Cells.Item(18, 12).Value = 900
$ws.Cells.Item(18, 13).Value = -152240.996
$ws.Cells.Item(18, 14).Value = -1238

$ws.Cells.Item(68, 8).Value = 1270.9012
$ws.Cells.Item(68, 10).Value = 1352.7213
$ws.Cells.Item(68, 12).Value = 4058.1639
$ws.Cells.Item(68, 14).Value = -5680.1639

$ws.Cells.Item(71, 8).Value = 1270.9012
$ws.Cells.Item(71, 10).Value = 1352.7213
$ws.Cells.Item(71, 12).Value = 12174.4917
$ws.Cells.Item(71, 14).Value = -20286.4917

$ws.Cells.Item(113, 8).Value = 2210.918
$ws.Cells.Item(113, 9).Value = 2995.8
$ws.Cells.Item(113, 10).Value = 715.9048
$ws.Cells.Item(113, 11).Value = 8987.400000000001
$ws.Cells.Item(113, 12).Value = 2147.7144
$ws.Cells.Item(113, 13).Value = -6817.400000000001
$ws.Cells.Item(113, 14).Value = -6487.7144

$ws.Cells.Item(140, 8).Value = 55920.95
$ws.Cells.Item(140, 9).Value = 167948.5
$ws.Cells.Item(140, 10).Value = 4215.923
$ws.Cells.Item(140, 11).Value = 503845.5
$ws.Cells.Item(140, 12).Value = 12647.769
$ws.Cells.Item(140, 13).Value = -498665.5
$ws.Cells.Item(140, 14).Value = -23007.769

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 1004051
$ws.Cells.Item(2, 9).Value = 1095323.9
$ws.Cells.Item(2, 10).Value = 50
$ws.Cells.Item(2, 11).Value = 1095323.9
$ws.Cells.Item(2, 12).Value = 50
$ws.Cells.Item(2, 13).Value = -1095210.9
$ws.Cells.Item(2, 14).Value = -276

$ws.Cells.Item(105, 8).Value = 41500
$ws.Cells.Item(105, 10).Value = 41500
$ws.Cells.Item(105, 12).Value = 41500
$ws.Cells.Item(105, 14).Value = -48488

$ws.Cells.Item(126, 8).Value = 8721.529
$ws.Cells.Item(126, 9).Value = 55006
$ws.Cells.Item(126, 10).Value = 2550.2666
$ws.Cells.Item(126, 11).Value = 165018
$ws.Cells.Item(126, 12).Value = 7650.7998
$ws.Cells.Item(126, 13).Value = -162548
$ws.Cells.Item(126, 14).Value = -12590.7998

$ws.Cells.Item(128, 8).Value = 40000
$ws.Cells.Item(128, 10).Value = 40000
$ws.Cells.Item(128, 12).Value = 40000
$ws.Cells.Item(128, 14).Value = -49960

$ws.Cells.Item(136, 8).Value = 29980
$ws.Cells.Item(136, 10).Value = 29980
$ws.Cells.Item(136, 12).Value = 89940
$ws.Cells.Item(136, 14).Value = -95040

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 4173.5806
$ws.Cells.Item(132, 9).Value = 3977.3076
$ws.Cells.Item(132, 10).Value = 4315.3335
$ws.Cells.Item(132, 11).Value = 11931.9228
$ws.Cells.Item(132, 12).Value = 12946.0005
$ws.Cells.Item(132, 13).Value = -9401.9228
$ws.Cells.Item(132, 14).Value = -18006.0005

$ws.Cells.Item(136, 8).Value = 1649.3429
$ws.Cells.Item(136, 9).Value = 1204.9
$ws.Cells.Item(136, 10).Value = 4316
$ws.Cells.Item(136, 11).Value = 3614.7
$ws.Cells.Item(136, 12).Value = 12948
$ws.Cells.Item(136, 13).Value = -1064.7
$ws.Cells.Item(136, 14).Value = -18048

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1226883.5
$ws.Cells.Item(126, 9).Value = 1401820.2
$ws.Cells.Item(126, 10).Value = 2326.6667
$ws.Cells.Item(126, 11).Value = 4205460.6
$ws.Cells.Item(126, 12).Value = 6980.000100000001
$ws.Cells.Item(126, 13).Value = -4202990.6
$ws.Cells.Item(126, 14).Value = -11920.0001

$ws.Cells.Item(132, 8).Value = 1612445.1
$ws.Cells.Item(132, 9).Value = 2289743.5
$ws.Cells.Item(132, 10).Value = 3861.375
$ws.Cells.Item(132, 11).Value = 6869230.5
$ws.Cells.Item(132, 12).Value = 11584.125
$ws.Cells.Item(132, 13).Value = -6866700.5
$ws.Cells.Item(132, 14).Value = -16644.125

Write-Output "Applied all cell updates."